# Applies the "Quick jump to specific component panel" slide edits:
#  - Curved Connector 10: tiny re-positioning/resizing + adj1 tweak
#    (off.x 357665->357666, ext.cx 99535->99534, adj1 -229668->-229670)
#  - TextBox 11: widened (ext.cx 2700304->3995805) and its text extended
#    from "Quick jump to specific component panel" to
#    "Quick jump to specific component panel or opens it if closed"
#    (split across two runs, matching how PowerPoint records a live edit)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Curved Connector 10 (id=11) ---
$connector = $s.Shapes.Item(7)
$connector.Left = 28.16267716535433      # EMU 357666
$connector.Width = 7.837323034645669     # EMU 99534
$connector.Adjustments.Item(1) = -2.2967 # val -229670

# --- TextBox 11 (id=12) ---
$textBox = $s.Shapes.Item(8)
$textBox.Width = 314.6303253606299       # EMU 3995805

$tr = $textBox.TextFrame.TextRange
$tr.Text = "Quick jump to specific component "
[void]$tr.InsertAfter("panel or opens it if closed")
